$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.629.19'
$ws.Range("E2").Value = '  +4.03%  '

$ws.Range("D3").Value = '1.798.78'
$ws.Range("E3").Value = '  +0.57%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.55'
$ws.Range("E5").Value = '  -0.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5314'
$ws.Range("E7").Value = '  -0.40%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3762'
$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07535'
$ws.Range("E9").Value = '  +0.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.56'
$ws.Range("E10").Value = '  -1.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.123'
$ws.Range("E11").Value = '  +1.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.20'
$ws.Range("E12").Value = '  +1.70%  '

$ws.Range("E13").Value = '  +0.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.204'
$ws.Range("E14").Value = '  +1.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.460'
$ws.Range("E15").Value = '  +5.89%  '

$ws.Range("D16").Value = '1.793.06'
$ws.Range("E16").Value = '  +0.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.56'
$ws.Range("E17").Value = '  +0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001069'
$ws.Range("E18").Value = '  +0.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06457'
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  +0.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.33'
$ws.Range("E21").Value = '  +2.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.930'
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").Value = '28.649.48'
$ws.Range("E23").Value = '  +4.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.21'
$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.094'
$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.96'
$ws.Range("E26").Value = '  +3.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.57'
$ws.Range("E27").Value = '  +0.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.398'
$ws.Range("E28").Value = '  +0.33%  '

$ws.Range("D29").Value = '2.003.71'
$ws.Range("E29").Value = '  +0.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.69'
$ws.Range("E30").Value = '  +1.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.133'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1022'
$ws.Range("E32").Value = '  -0.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.728'
$ws.Range("E33").Value = '  +1.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.665'
$ws.Range("E34").Value = '  +1.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2304'
$ws.Range("E35").Value = '  +11.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06560'
$ws.Range("E36").Value = '  +9.18%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02324'
$ws.Range("E37").Value = '  +2.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.817'
$ws.Range("E38").Value = '  +2.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.081'
$ws.Range("E39").Value = '  +2.48%  '

$ws.Range("E40").Value = '  +1.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6320'
$ws.Range("E41").Value = '  +2.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.200'
$ws.Range("E42").Value = '  +4.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.394'
$ws.Range("E44").Value = '  -1.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.58'
$ws.Range("E45").Value = '  +1.15%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5928'
$ws.Range("E46").Value = '  +1.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.666'
$ws.Range("E47").Value = '  +0.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.86'
$ws.Range("E48").Value = '  +3.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.982'
$ws.Range("E49").Value = '  +3.91%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.172'
$ws.Range("E50").Value = '  +3.98%  '

$ws.Range("E51").Value = '  +2.88%  '
